$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F holds another "test date" column, same date-number style as
# the existing C:E columns (style comes from E1).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = 41382

# New result recorded for row 10 ("OK") and the row grows a bit taller to
# fit it.
$ws.Range("F10").Value = "OK"
$ws.Rows(10).RowHeight = 16.5

# The selection moves to the newly-entered cell (this also drops the old
# scrolled-down view position).
$ws.Range("F10").Select()
